$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells whose new values look like plain numbers stay as text
# (preserving formatting such as trailing zeros), matching the source data.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = "58.286.41"
$ws.Range("E2").Value = "  -2.88%  "
$ws.Range("D3").Value = "2.965.29"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "557.80"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").Value = "131.73"
$ws.Range("E6").Value = "  +7.83%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").Value = "0.519"
$ws.Range("E8").Value = "  +4.82%  "
$ws.Range("D9").Value = "2.959.55"
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").Value = "0.130"
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("D11").Value = "4.85"
$ws.Range("E11").Value = "  -4.21%  "
$ws.Range("D12").Value = "0.450"
$ws.Range("E12").Value = "  +4.07%  "
$ws.Range("E13").Value = "  +2.76%  "
$ws.Range("D14").Value = "32.97"
$ws.Range("E14").Value = "  +2.66%  "
$ws.Range("E15").Value = "  +2.43%  "
$ws.Range("D16").Value = "3.468.65"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("D17").Value = "6.81"
$ws.Range("E17").Value = "  +11.72%  "
$ws.Range("D18").Value = "2.972.51"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").Value = "58.542.00"
$ws.Range("E19").Value = "  -2.62%  "
$ws.Range("D20").Value = "421.64"
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("D21").Value = "13.19"
$ws.Range("E21").Value = "  +2.36%  "
$ws.Range("D22").Value = "0.686"
$ws.Range("D23").Value = "7.01"
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("D24").Value = "13.01"
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("D25").Value = "79.75"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").Value = "7.59"
$ws.Range("E29").Value = "  +7.57%  "
$ws.Range("D30").Value = "2.01"
$ws.Range("E30").Value = "  +8.61%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "6.20"
$ws.Range("E31").Value = "  +2.67%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.105"
$ws.Range("E32").Value = "  +14.73%  "
$ws.Range("D33").Value = "25.22"
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("B34").Value = "Stacks"
$ws.Range("C34").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D34").Value = "2.13"
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "5.67"
$ws.Range("E35").Value = "  +3.25%  "
$ws.Range("D36").Value = "0.948"
$ws.Range("E36").Value = "  +1.23%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "48.63"
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0691"
$ws.Range("E38").Value = "  +8.19%  "
$ws.Range("D39").Value = "8.44"
$ws.Range("E39").Value = "  +8.43%  "
$ws.Range("E40").Value = "  +10.70%  "
$ws.Range("D41").Value = "0.109"
$ws.Range("E41").Value = "  +1.66%  "
$ws.Range("D42").Value = "0.0352"
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("D43").Value = "380.31"
$ws.Range("E43").Value = "  +1.95%  "
$ws.Range("D44").Value = "2.661.66"
$ws.Range("E44").Value = "  +2.18%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "0.241"
$ws.Range("E46").Value = "  +3.82%  "
$ws.Range("D47").Value = "121.36"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("E48").Value = "  +4.27%  "
$ws.Range("D49").Value = "1.99"
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("D50").Value = "23.52"
$ws.Range("E50").Value = "  +2.48%  "
$ws.Range("D51").Value = "2.01"
$ws.Range("E51").Value = "  +2.81%  "
